{"js": "// Office.js (Word JavaScript API) edit script.\n//\n// Target change: the bullet paragraph\n//   \"Messages: SIDs URNs Case Classes Statements. Statement Data Pattern\n//    Matching. State Flows: Reactive Events Messages.\"\n// loses the stray paragraph-mark underline formatting it carried, and is\n// followed by seven new bullets (same numbered-list, numId 3) ending with\n// \"Sample Workflow: ToDo\".\n\n// Escape text for safe embedding inside XML.\nfunction xmlEscape(s) {\n  return s\n    .replace(/&/g, \"&amp;\")\n    .replace(/</g, \"&lt;\")\n    .replace(/>/g, \"&gt;\")\n    .replace(/\"/g, \"&quot;\")\n    .replace(/'/g, \"&apos;\");\n}\n\n// The anchor paragraph's exact current text (used to locate it) and the\n// full ordered list of bullet texts that should exist in its place\n// afterwards (first entry = the anchor paragraph itself, with its\n// paragraph-mark underline override removed; the rest are brand-new\n// bullets appended after it).\nconst anchorText =\n  \"Messages: SIDs URNs Case Classes Statements. Statement Data Pattern Matching. State Flows: Reactive Events Messages.\";\n\nconst newBulletTexts = [\n  \"Protocol: SIDs URNs Resources. Endpoints: Case Classes Events Signatures, Statement Data Pattern Matching Events.\",\n  \"Core Model Upper Resources (DCI Context / Facets: Metaclass, Class, etc. as Resource, root navigation Context Resource).\",\n  \"Protocol: GET URN Case Classes / Statement Data Aggregated Events Messages Statements.\",\n  \"Protocol: Browse Messages Events Statements. Build Context State Flows.\",\n  \"Protocol: POST URN Navigation Context built Case Class Statement Data Events.\",\n  \"Protocol: POST Subsequent entailed Context Browsing / Events Transforms.\",\n  \"Sample Workflow: ToDo\",\n];\n\nconst allTexts = [anchorText, ...newBulletTexts];\n\n// Build one <w:p> per bullet: same list (numId 3, ilvl 0) and indent the\n// anchor paragraph already used, but with a clean paragraph mark (no\n// <w:rPr><w:u w:val=\"none\"/></w:rPr> override) and a single plain run.\nfunction paragraphXml(text) {\n  return (\n    \"<w:p><w:pPr><w:numPr><w:ilvl w:val=\\\"0\\\"/><w:numId w:val=\\\"3\\\"/></w:numPr>\" +\n    \"<w:ind w:left=\\\"600\\\" w:hanging=\\\"360\\\"/></w:pPr>\" +\n    \"<w:r><w:rPr><w:rtl w:val=\\\"0\\\"/></w:rPr><w:t xml:space=\\\"preserve\\\">\" +\n    xmlEscape(text) +\n    \"</w:t></w:r></w:p>\"\n  );\n}\n\nconst bodyXml = allTexts.map(paragraphXml).join(\"\");\n\nconst flatOpcXml =\n  '<?xml version=\"1.0\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  \"<w:body>\" +\n  bodyXml +\n  \"</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>\";\n\n// Locate the anchor paragraph via a body search on its exact text.\nconst body = context.document.body;\nconst results = body.search(anchorText, { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Could not find the target paragraph to update.\");\n}\n\nconst anchorParagraph = results.items[0].paragraphs.getFirst();\nconst wholeRange = anchorParagraph.getRange(\"Whole\");\n\n// Replace the single anchor paragraph with the full run of (cleaned\n// anchor + 7 new) bullet paragraphs in one shot.\nwholeRange.insertOoxml(flatOpcXml, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Word COM interop edit script.\n#\n# Target change: the bullet paragraph\n#   \"Messages: SIDs URNs Case Classes Statements. Statement Data Pattern\n#    Matching. State Flows: Reactive Events Messages.\"\n# loses the stray paragraph-mark underline formatting it carried, and is\n# followed by seven new bullets (same numbered-list, numId 3) ending with\n# \"Sample Workflow: ToDo\".\n\n$d = $word.ActiveDocument\n\n# The anchor paragraph's exact current text (used to locate it) and the\n# full ordered list of bullet texts that should exist in its place\n# afterwards (first entry = the anchor paragraph itself, with its\n# paragraph-mark underline override removed; the rest are brand-new\n# bullets appended after it).\n$anchorText = \"Messages: SIDs URNs Case Classes Statements. Statement Data Pattern Matching. State Flows: Reactive Events Messages.\"\n\n$newBulletTexts = @(\n  \"Protocol: SIDs URNs Resources. Endpoints: Case Classes Events Signatures, Statement Data Pattern Matching Events.\",\n  \"Core Model Upper Resources (DCI Context / Facets: Metaclass, Class, etc. as Resource, root navigation Context Resource).\",\n  \"Protocol: GET URN Case Classes / Statement Data Aggregated Events Messages Statements.\",\n  \"Protocol: Browse Messages Events Statements. Build Context State Flows.\",\n  \"Protocol: POST URN Navigation Context built Case Class Statement Data Events.\",\n  \"Protocol: POST Subsequent entailed Context Browsing / Events Transforms.\",\n  \"Sample Workflow: ToDo\"\n)\n\n$allTexts = @($anchorText) + $newBulletTexts\n\n# Build one <w:p> per bullet: same list (numId 3, ilvl 0) and indent the\n# anchor paragraph already used, but with a clean paragraph mark (no\n# <w:rPr><w:u w:val=\"none\"/></w:rPr> override) and a single plain run.\n$bodyXml = \"\"\nforeach ($t in $allTexts) {\n  $escaped = $t.Replace(\"&\", \"&amp;\").Replace(\"<\", \"&lt;\").Replace(\">\", \"&gt;\").Replace(\"'\", \"&apos;\").Replace('\"', \"&quot;\")\n  $bodyXml += '<w:p><w:pPr><w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"3\"/></w:numPr><w:ind w:left=\"600\" w:hanging=\"360\"/></w:pPr><w:r><w:rPr><w:rtl w:val=\"0\"/></w:rPr><w:t xml:space=\"preserve\">' + $escaped + '</w:t></w:r></w:p>'\n}\n\n$flatOpcXml = '<?xml version=\"1.0\" standalone=\"yes\"?><pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body>' + $bodyXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n\n# Locate the anchor paragraph via Find, then replace just that one\n# paragraph with the full run of (cleaned anchor + 7 new) bullets.\n$range = $d.Content\n$found = $range.Find.Execute($anchorText)\nif (-not $found) {\n  throw \"Could not find the target paragraph to update.\"\n}\n\n$para = $range.Paragraphs(1)\n$para.Range.InsertXML($flatOpcXml)\n"}
